$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '30.451.38'
Set-TextValue $ws.Range('E2') '  +0.80%  '
Set-TextValue $ws.Range('D3') '2.106.39'
Set-TextValue $ws.Range('E3') '  +1.73%  '
Set-TextValue $ws.Range('D4') '1.008'
Set-TextValue $ws.Range('D5') '334.81'
Set-TextValue $ws.Range('E5') '  +2.51%  '
Set-TextValue $ws.Range('D6') '1.006'
Set-TextValue $ws.Range('E6') '  +0.83%  '
Set-TextValue $ws.Range('D7') '0.5228'
Set-TextValue $ws.Range('E7') '  +0.94%  '
Set-TextValue $ws.Range('D8') '0.4543'
Set-TextValue $ws.Range('E8') '  +5.18%  '
Set-TextValue $ws.Range('D9') '52.94'
Set-TextValue $ws.Range('E9') '  +15.76%  '
Set-TextValue $ws.Range('D10') '0.08922'
Set-TextValue $ws.Range('E10') '  +0.14%  '
Set-TextValue $ws.Range('D11') '1.175'
Set-TextValue $ws.Range('E11') '  +2.06%  '
Set-TextValue $ws.Range('D12') '24.27'
Set-TextValue $ws.Range('E12') '  +0.49%  '
Set-TextValue $ws.Range('D13') '2.117.12'
Set-TextValue $ws.Range('E13') '  +2.13%  '
Set-TextValue $ws.Range('D14') '6.854'
Set-TextValue $ws.Range('E14') '  +3.22%  '
Set-TextValue $ws.Range('D15') '8.026'
Set-TextValue $ws.Range('E15') '  +5.08%  '
Set-TextValue $ws.Range('D16') '96.39'
Set-TextValue $ws.Range('E16') '  +1.73%  '
Set-TextValue $ws.Range('D17') '0.00001145'
Set-TextValue $ws.Range('E17') '  +2.30%  '
Set-TextValue $ws.Range('D18') '1.007'
Set-TextValue $ws.Range('E18') '  +0.75%  '
Set-TextValue $ws.Range('D19') '0.06658'
Set-TextValue $ws.Range('E19') '  +0.86%  '
Set-TextValue $ws.Range('D20') '19.21'
Set-TextValue $ws.Range('E20') '  +2.73%  '
Set-TextValue $ws.Range('E21') '  +0.77%  '
Set-TextValue $ws.Range('D22') '6.353'
Set-TextValue $ws.Range('E22') '  +2.49%  '
Set-TextValue $ws.Range('D23') '30.503.28'
Set-TextValue $ws.Range('E23') '  +0.80%  '
Set-TextValue $ws.Range('D24') '12.40'
Set-TextValue $ws.Range('E24') '  +1.63%  '
Set-TextValue $ws.Range('D25') '2.374'
Set-TextValue $ws.Range('E25') '  +4.10%  '
Set-TextValue $ws.Range('D26') '2.366.93'
Set-TextValue $ws.Range('E26') '  +2.24%  '
Set-TextValue $ws.Range('D27') '22.26'
Set-TextValue $ws.Range('E27') '  +0.57%  '
Set-TextValue $ws.Range('D28') '163.60'
Set-TextValue $ws.Range('E28') '  +1.36%  '
Set-TextValue $ws.Range('D29') '2.543'
Set-TextValue $ws.Range('E29') '  +0.91%  '
Set-TextValue $ws.Range('D30') '132.81'
Set-TextValue $ws.Range('E30') '  +1.67%  '
Set-TextValue $ws.Range('D31') '1.221'
Set-TextValue $ws.Range('E31') '  +2.98%  '
Set-TextValue $ws.Range('D32') '0.1071'
Set-TextValue $ws.Range('E32') '  +0.64%  '
Set-TextValue $ws.Range('D33') '1.658'
Set-TextValue $ws.Range('E33') '  +2.40%  '
Set-TextValue $ws.Range('D34') '6.327'
Set-TextValue $ws.Range('E34') '  +4.11%  '
Set-TextValue $ws.Range('D35') '3.947'
Set-TextValue $ws.Range('E35') '  +3.32%  '
Set-TextValue $ws.Range('E36') '  +8.69%  '
Set-TextValue $ws.Range('D37') '0.02580'
Set-TextValue $ws.Range('E37') '  +0.70%  '
Set-TextValue $ws.Range('D38') '5.683'
Set-TextValue $ws.Range('E38') '  +4.94%  '
Set-TextValue $ws.Range('D39') '0.06833'
Set-TextValue $ws.Range('E39') '  +3.33%  '
Set-TextValue $ws.Range('D40') '0.2301'
Set-TextValue $ws.Range('E40') '  +2.98%  '
Set-TextValue $ws.Range('D41') '12.63'
Set-TextValue $ws.Range('E41') '  +0.88%  '
Set-TextValue $ws.Range('D42') '0.6881'
Set-TextValue $ws.Range('E42') '  +1.12%  '
Set-TextValue $ws.Range('D43') '1.252'
Set-TextValue $ws.Range('E43') '  +0.91%  '
Set-TextValue $ws.Range('B44') 'Frax'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D44') '1.006'
Set-TextValue $ws.Range('E44') '  +0.88%  '
Set-TextValue $ws.Range('B45') 'NEARProtocol'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D45') '2.335'
Set-TextValue $ws.Range('E45') '  +6.41%  '
Set-TextValue $ws.Range('D46') '14.08'
Set-TextValue $ws.Range('E46') '  +1.44%  '
Set-TextValue $ws.Range('D47') '0.6371'
Set-TextValue $ws.Range('E47') '  +0.55%  '
Set-TextValue $ws.Range('D48') '3.667'
Set-TextValue $ws.Range('E48') '  +2.05%  '
Set-TextValue $ws.Range('D49') '1.248'
Set-TextValue $ws.Range('E49') '  +1.42%  '
Set-TextValue $ws.Range('D50') '0.3422'
Set-TextValue $ws.Range('E50') '  +25.97%  '
Set-TextValue $ws.Range('D51') '83.30'
Set-TextValue $ws.Range('E51') '  +2.84%  '
